$d = $word.ActiveDocument

# The document ends with an automatic "_GoBack" bookmark marking the position
# of the last edit. Remove it now; we will re-create it at the correct final
# position once all of the new content below has been inserted.
$bmGoBack = $d.Bookmarks("_GoBack")
$bmGoBack.Delete()

# --- 1) Append a new sentence to the end of the last paragraph ----------
# (the sentence ending "...kulturelle og nationale oplevelser." gets a
# trailing sentence appended to it)
$lastPara = $d.Paragraphs.Last
$insertPos = $lastPara.Range.End - 1   # just before the paragraph mark
$r1 = $d.Range($insertPos, $insertPos)
$r1.InsertAfter(" Det er dog individuelt, hvad en turist mener er interessant.")
$r1.LanguageID = "da-DK"

# --- 2) Insert a brand new paragraph right after that one ----------------
$r2 = $d.Range($r1.End, $r1.End)
$r2.InsertParagraphAfter()

$newParaStart = $r2.End + 1
$r3 = $d.Range($newParaStart, $newParaStart)
$r3.InsertAfter("De sidste 20% fra vores sp" + [char]0x00F8 + "rgeskema-unders" + [char]0x00F8 + "gelse, som valgte den hurtigste rute, vil have en rute mellem de attraktioner de " + [char]0x00F8 + "nsker at se, som enten er den korteste i afstand, eller i tid. Den interessante rute vil derved v" + [char]0x00E6 + "re en udvidelse af dette, da det vil foresl" + [char]0x00E5 + " en muligvis l" + [char]0x00E6 + "ngere rute, som indeholder flere attraktioner eller oplevelser undervejs.")
$r3.LanguageID = "da-DK"

# --- 3) Re-create the "_GoBack" bookmark at the new end of the document --
# Adding a bookmark exactly at the end-of-story position is unreliable in
# this runtime, so pad with one sacrificial character, add the bookmark
# just before it, then delete the padding again.
$endPos = $d.Content.End
$padRange = $d.Range($endPos - 1, $endPos - 1)
$padRange.InsertAfter("X")

$bmPos = $endPos - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$padDelRange = $d.Range($bmPos, $bmPos + 1)
$padDelRange.Delete()
